$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.276052666666667
$ws.Range("H2").Value = 6.828158
$ws.Range("I2").Value = 0.005247614157263819
$ws.Range("J2").Value = 0.005247614157263819
$ws.Range("M2").Value = 3.135398666666667
$ws.Range("N2").Value = 9.406196000000001
$ws.Range("O2").Value = 0.1723049126704688
$ws.Range("P2").Value = 0.1723049126704688
$ws.Range("Q2").Value = 7.136332496329779
$ws.Range("R2").Value = 64.22699246696801
$ws.Range("S2").Value = 0.0009041896990956581
$ws.Range("T2").Value = 0.0009041896990956582

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.276052666666667
$ws.Range("H3").Value = 6.828158
$ws.Range("I3").Value = 0.005247614157263819
$ws.Range("J3").Value = 0.005247614157263819
$ws.Range("O3").Value = 0.1733096678828815
$ws.Range("P3").Value = 0.1733096678828815
$ws.Range("Q3").Value = 7.177946325918667
$ws.Range("R3").Value = 64.601516933268
$ws.Range("S3").Value = 0.0009094622667728993
$ws.Range("T3").Value = 0.0009094622667728994

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.276052666666667
$ws.Range("H4").Value = 6.828158
$ws.Range("I4").Value = 0.005247614157263819
$ws.Range("J4").Value = 0.005247614157263819
$ws.Range("M4").Value = 0.4900660000000001
$ws.Range("N4").Value = 1.470198
$ws.Range("O4").Value = 0.02693143306797965
$ws.Range("P4").Value = 0.02693143306797965
$ws.Range("Q4").Value = 1.115416026142667
$ws.Range("R4").Value = 10.038744235284
$ws.Range("S4").Value = 0.000141325769442933
$ws.Range("T4").Value = 0.000141325769442933

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.276052666666667
$ws.Range("H5").Value = 6.828158
$ws.Range("I5").Value = 0.005247614157263819
$ws.Range("J5").Value = 0.005247614157263819
$ws.Range("M5").Value = 11.417657
$ws.Range("N5").Value = 34.252971
$ws.Range("O5").Value = 0.62745398637867
$ws.Range("P5").Value = 0.6274539863786701
$ws.Range("Q5").Value = 25.98718866193533
$ws.Range("R5").Value = 233.884697957418
$ws.Range("S5").Value = 0.003292636421952328
$ws.Range("T5").Value = 0.003292636421952329

$ws.Range("I6").Value = 0.1062533062835484
$ws.Range("J6").Value = 0.1062533062835484
$ws.Range("M6").Value = 3.135398666666667
$ws.Range("N6").Value = 9.406196000000001
$ws.Range("O6").Value = 0.1723049126704688
$ws.Range("P6").Value = 0.1723049126704688
$ws.Range("Q6").Value = 144.4959365818036
$ws.Range("R6").Value = 1300.463429236232
$ws.Range("S6").Value = 0.01830796666013538
$ws.Range("T6").Value = 0.01830796666013538

$ws.Range("I7").Value = 0.1062533062835484
$ws.Range("J7").Value = 0.1062533062835484
$ws.Range("O7").Value = 0.1733096678828815
$ws.Range("P7").Value = 0.1733096678828815
$ws.Range("S7").Value = 0.01841472522345985
$ws.Range("T7").Value = 0.01841472522345985

$ws.Range("I8").Value = 0.1062533062835484
$ws.Range("J8").Value = 0.1062533062835484
$ws.Range("M8").Value = 0.4900660000000001
$ws.Range("N8").Value = 1.470198
$ws.Range("O8").Value = 0.02693143306797965
$ws.Range("P8").Value = 0.02693143306797965
$ws.Range("Q8").Value = 22.58486182625734
$ws.Range("R8").Value = 203.263756436316
$ws.Range("S8").Value = 0.002861553806426925
$ws.Range("T8").Value = 0.002861553806426925

$ws.Range("I9").Value = 0.1062533062835484
$ws.Range("J9").Value = 0.1062533062835484
$ws.Range("M9").Value = 11.417657
$ws.Range("N9").Value = 34.252971
$ws.Range("O9").Value = 0.62745398637867
$ws.Range("P9").Value = 0.6274539863786701
$ws.Range("Q9").Value = 526.1866885778647
$ws.Range("R9").Value = 4735.680197200782
$ws.Range("S9").Value = 0.06666906059352622
$ws.Range("T9").Value = 0.06666906059352624

$ws.Range("G10").Value = 41.187613
$ws.Range("H10").Value = 123.562839
$ws.Range("I10").Value = 0.09496120377532416
$ws.Range("J10").Value = 0.09496120377532417
$ws.Range("M10").Value = 3.135398666666667
$ws.Range("N10").Value = 9.406196000000001
$ws.Range("O10").Value = 0.1723049126704688
$ws.Range("P10").Value = 0.1723049126704688
$ws.Range("Q10").Value = 129.1395868833827
$ws.Range("R10").Value = 1162.256281950444
$ws.Range("S10").Value = 0.01636228192358983
$ws.Range("T10").Value = 0.01636228192358983

$ws.Range("G11").Value = 41.187613
$ws.Range("H11").Value = 123.562839
$ws.Range("I11").Value = 0.09496120377532416
$ws.Range("J11").Value = 0.09496120377532417
$ws.Range("O11").Value = 0.1733096678828815
$ws.Range("P11").Value = 0.1733096678828815
$ws.Range("Q11").Value = 129.892633741066
$ws.Range("R11").Value = 1169.033703669594
$ws.Range("S11").Value = 0.01645769468806006
$ws.Range("T11").Value = 0.01645769468806006

$ws.Range("G12").Value = 41.187613
$ws.Range("H12").Value = 123.562839
$ws.Range("I12").Value = 0.09496120377532416
$ws.Range("J12").Value = 0.09496120377532417
$ws.Range("M12").Value = 0.4900660000000001
$ws.Range("N12").Value = 1.470198
$ws.Range("O12").Value = 0.02693143306797965
$ws.Range("P12").Value = 0.02693143306797965
$ws.Range("Q12").Value = 20.184648752458
$ws.Range("R12").Value = 181.661838772122
$ws.Range("S12").Value = 0.002557441303529919
$ws.Range("T12").Value = 0.00255744130352992

$ws.Range("G13").Value = 41.187613
$ws.Range("H13").Value = 123.562839
$ws.Range("I13").Value = 0.09496120377532416
$ws.Range("J13").Value = 0.09496120377532417
$ws.Range("M13").Value = 11.417657
$ws.Range("N13").Value = 34.252971
$ws.Range("O13").Value = 0.62745398637867
$ws.Range("P13").Value = 0.6274539863786701
$ws.Range("Q13").Value = 470.266037882741
$ws.Range("R13").Value = 4232.394340944669
$ws.Range("S13").Value = 0.05958378586014435
$ws.Range("T13").Value = 0.05958378586014437

$ws.Range("G14").Value = 344.1819356666667
$ws.Range("H14").Value = 1032.545807
$ws.Range("I14").Value = 0.7935378757838636
$ws.Range("J14").Value = 0.7935378757838637
$ws.Range("M14").Value = 3.135398666666667
$ws.Range("N14").Value = 9.406196000000001
$ws.Range("O14").Value = 0.1723049126704688
$ws.Range("P14").Value = 0.1723049126704688
$ws.Range("Q14").Value = 1079.147582180019
$ws.Range("R14").Value = 9712.328239620174
$ws.Range("S14").Value = 0.136730474387648
$ws.Range("T14").Value = 0.136730474387648

$ws.Range("G15").Value = 344.1819356666667
$ws.Range("H15").Value = 1032.545807
$ws.Range("I15").Value = 0.7935378757838636
$ws.Range("J15").Value = 0.7935378757838637
$ws.Range("O15").Value = 0.1733096678828815
$ws.Range("P15").Value = 0.1733096678828815
$ws.Range("Q15").Value = 1085.440375237125
$ws.Range("R15").Value = 9768.963377134121
$ws.Range("S15").Value = 0.1375277857045887
$ws.Range("T15").Value = 0.1375277857045887

$ws.Range("G16").Value = 344.1819356666667
$ws.Range("H16").Value = 1032.545807
$ws.Range("I16").Value = 0.7935378757838636
$ws.Range("J16").Value = 0.7935378757838637
$ws.Range("M16").Value = 0.4900660000000001
$ws.Range("N16").Value = 1.470198
$ws.Range("O16").Value = 0.02693143306797965
$ws.Range("P16").Value = 0.02693143306797965
$ws.Range("Q16").Value = 168.6718644844207
$ws.Range("R16").Value = 1518.046780359786
$ws.Range("S16").Value = 0.02137111218857987
$ws.Range("T16").Value = 0.02137111218857988

$ws.Range("G17").Value = 344.1819356666667
$ws.Range("H17").Value = 1032.545807
$ws.Range("I17").Value = 0.7935378757838636
$ws.Range("J17").Value = 0.7935378757838637
$ws.Range("M17").Value = 11.417657
$ws.Range("N17").Value = 34.252971
$ws.Range("O17").Value = 0.62745398637867
$ws.Range("P17").Value = 0.6274539863786701
$ws.Range("Q17").Value = 3929.751287038066
$ws.Range("R17").Value = 35367.7615833426
$ws.Range("S17").Value = 0.4979085035030471
$ws.Range("T17").Value = 0.4979085035030472
